$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 25 / Row 26: swap A/B/Q/R values ---
$a25 = $ws.Range("A25").Value()
$b25 = $ws.Range("B25").Value()
$q25 = $ws.Range("Q25").Value()
$r25 = $ws.Range("R25").Value()

$a26 = $ws.Range("A26").Value()
$b26 = $ws.Range("B26").Value()
$q26 = $ws.Range("Q26").Value()
$r26 = $ws.Range("R26").Value()

$ws.Range("A25").Value() = $a26
$ws.Range("Q25").Value() = $q26
$ws.Range("R25").Value() = $r26

$ws.Range("A26").Value() = $a25
$ws.Range("Q26").Value() = $q25
$ws.Range("R26").Value() = $r25

# B25 and B26 both become 89571
$ws.Range("B25").Value() = 89571
$ws.Range("B26").Value() = 89571

# --- Rows 27 / 30 / 31: rotate species-related data (E,F,G,H,Q,R) and A ---
# Capture "before" state of the three rows first.
$a27 = $ws.Range("A27").Value()
$e27 = $ws.Range("E27").Value()
$f27 = $ws.Range("F27").Value()
$g27 = $ws.Range("G27").Value()
$h27 = $ws.Range("H27").Value()
$q27 = $ws.Range("Q27").Value()
$r27 = $ws.Range("R27").Value()

$a30 = $ws.Range("A30").Value()
$e30 = $ws.Range("E30").Value()
$f30 = $ws.Range("F30").Value()
$g30 = $ws.Range("G30").Value()
$h30 = $ws.Range("H30").Value()
$q30 = $ws.Range("Q30").Value()
$r30 = $ws.Range("R30").Value()

$a31 = $ws.Range("A31").Value()
$e31 = $ws.Range("E31").Value()
$f31 = $ws.Range("F31").Value()
$g31 = $ws.Range("G31").Value()
$h31 = $ws.Range("H31").Value()
$q31 = $ws.Range("Q31").Value()
$r31 = $ws.Range("R31").Value()

# Row 27 gets row 31's former species/location data
$ws.Range("A27").Value() = $a31
$ws.Range("B27").Value() = 56430
$ws.Range("E27").Value() = $e31
$ws.Range("F27").Value() = $f31
$ws.Range("G27").Value() = $g31
$ws.Range("H27").Value() = $h31
$ws.Range("Q27").Value() = $q31
$ws.Range("R27").Value() = $r31

# Row 30 gets row 27's former species/location data
$ws.Range("A30").Value() = $a27
$ws.Range("B30").Value() = 77650
$ws.Range("E30").Value() = $e27
$ws.Range("F30").Value() = $f27
$ws.Range("G30").Value() = $g27
$ws.Range("H30").Value() = $h27
$ws.Range("Q30").Value() = $q27
$ws.Range("R30").Value() = $r27

# Row 31 gets row 30's former species/location data
$ws.Range("A31").Value() = $a30
$ws.Range("B31").Value() = 56446
$ws.Range("E31").Value() = $e30
$ws.Range("F31").Value() = $f30
$ws.Range("G31").Value() = $g30
$ws.Range("H31").Value() = $h30
$ws.Range("Q31").Value() = $q30
$ws.Range("R31").Value() = $r30

# Row 27 gains the blank K/L/M/N cells that row 31 used to have
$ws.Range("I27").Copy($ws.Range("K27"))
$ws.Range("I27").Copy($ws.Range("L27"))
$ws.Range("I27").Copy($ws.Range("M27"))
$ws.Range("I27").Copy($ws.Range("N27"))

# Row 27 gains the public-comment note that used to be on row 31
$ws.Range("AC27").Value() = "ringhack äldre"

# Row 31 loses its former blank K/L/M/N cells and its public comment
$ws.Range("K31").ClearContents()
$ws.Range("L31").ClearContents()
$ws.Range("M31").ClearContents()
$ws.Range("N31").ClearContents()
$ws.Range("AC31").ClearContents()

# --- Row 28 / Row 29: simple B value updates ---
$ws.Range("B28").Value() = 90113
$ws.Range("B29").Value() = 90235
